# Fix Student Excel Sheet
# - Rename Sheet1 -> "All Students", add a new "Batch Name" column in front
# - Add one new sheet per batch ("1234-1234", "2025-2026", "2025-2027") with
#   that batch's single row of upload data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the original sheet and rebuild its contents with the new
#    "Batch Name" column + merged "Upload Date & Time" column.
# ---------------------------------------------------------------------
$allStudents = $wb.Worksheets.Item(1)
$allStudents.Name = "All Students"

# Wipe everything (data + the bold/border header style) and start clean.
$allStudents.Cells.Clear()

$allStudents.Cells.Item(1,1).Value = "Batch Name"
$allStudents.Cells.Item(1,2).Value = "ER Number"
$allStudents.Cells.Item(1,3).Value = "Student Name"
$allStudents.Cells.Item(1,4).Value = "Upload Date & Time"

$batches = @(
    @{ Name = "1234-1234"; ER = "92310133004"; Student = "Bhargav_1"; When = "2025-09-07 12:17:35" },
    @{ Name = "2025-2026"; ER = "92310133004"; Student = "Bhargav_1"; When = "2025-09-07 12:17:42" },
    @{ Name = "2025-2027"; ER = "92310133004"; Student = "Bhargav_1"; When = "2025-09-07 12:19:09" }
)

$row = 2
foreach ($b in $batches) {
    $allStudents.Cells.Item($row,1).Value = $b.Name

    # ER numbers are digit strings that must stay text, not become numbers.
    $erCell = $allStudents.Cells.Item($row,2)
    $erCell.NumberFormat = "@"
    $erCell.Value = $b.ER
    $erCell.Style = "Normal"

    $allStudents.Cells.Item($row,3).Value = $b.Student
    $allStudents.Cells.Item($row,4).Value = $b.When
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 2. Add one worksheet per batch, each holding that batch's single row.
# ---------------------------------------------------------------------
$prevSheet = $allStudents
foreach ($b in $batches) {
    $newSheet = $wb.Worksheets.Add($null, $prevSheet)
    $newSheet.Name = $b.Name

    $newSheet.Cells.Item(1,1).Value = "ER Number"
    $newSheet.Cells.Item(1,2).Value = "Student Name"
    $newSheet.Cells.Item(1,3).Value = "Upload Date & Time"

    $erCell2 = $newSheet.Cells.Item(2,1)
    $erCell2.NumberFormat = "@"
    $erCell2.Value = $b.ER
    $erCell2.Style = "Normal"

    $newSheet.Cells.Item(2,2).Value = $b.Student
    $newSheet.Cells.Item(2,3).Value = $b.When

    $prevSheet = $newSheet
}
